$d = $word.ActiveDocument

# Locate the target paragraph: the "Resumen" body paragraph that starts with
# the quoted "Climbing Journey" sentence.
$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith('"Climbing Journ')) {
        $targetPara = $p
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not locate target paragraph"
}

$start = $targetPara.Range.Start
$end = $targetPara.Range.End

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512"><pkg:xmlData><Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships"><Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/></Relationships></pkg:xmlData></pkg:part><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00230832" w:rsidRPr="00230832" w:rsidRDefault="00A307AD" w:rsidP="00230832"><w:pPr><w:pStyle w:val="APASEPTIMA"/></w:pPr><w:r w:rsidRPr="00A307AD"><w:t>"Climbing Journ</w:t></w:r><w:r><w:t>ey" será un juego sobre escalada el cual presentara escenarios en los cuales se tendrá que escalar</w:t></w:r><w:r w:rsidRPr="00A307AD"><w:t xml:space="preserve"> en roca (Rock Climbing) o en hielo (Ice Climbing). El juego se desarrollará util</w:t></w:r><w:r><w:t xml:space="preserve">izando el motor gráfico </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Unity</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r w:rsidRPr="00A307AD"><w:t>Será un juego gratuito disponible en la plataforma itch.io, con opción a donar.</w:t></w:r><w:r><w:t xml:space="preserve"> El objetivo va a consistir en e</w:t></w:r><w:r><w:t>scalar</w:t></w:r><w:r><w:t xml:space="preserve"> hasta la cima</w:t></w:r><w:r><w:t>, con un</w:t></w:r><w:r><w:t xml:space="preserve">a dificultad Normal o en Hielo. </w:t></w:r><w:r><w:t xml:space="preserve">El juego se </w:t></w:r><w:r><w:t>desarrollará</w:t></w:r><w:r><w:t xml:space="preserve"> de manera individual. Al inicio de cada etapa, el jugador podrá realizar un relevamiento del “Mapa”, para inte</w:t></w:r><w:r><w:t>riorizarse de las dificultades. El juego estará regido por o</w:t></w:r><w:r><w:t xml:space="preserve">bjetivos y por </w:t></w:r><w:r><w:t>puntaje</w:t></w:r><w:r><w:t xml:space="preserve">. AL alcanzar la meta el jugador formará parte de un ranking, y va adquirir mejoras adicionales como premio (herramientas, calzados, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r><w:r><w:t>.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$replaceRange = $d.Range($start, $end)
$null = $replaceRange.InsertXML($xml)
